# "perbaiki laporan Lipa 14"
# - update the sign-off date and add a trailing space after "Panitera,"
# - bump the Pagu Awal / Sisa figures on row 9 (70,000,000 -> 90,000,000 ; 69,900,000 -> 89,900,000)
# - enlarge the font (11 -> 12 pt) for the summary table (row 9) and the
#   signature block (rows 12-19), and turn word-wrap on for row 9

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text updates -------------------------------------------------------
$ws.Range("H12").Value = "Ternate , 02 Agustus 2023"
$ws.Range("H13").Value = "Panitera, "

# --- Numeric updates ------------------------------------------------------
$ws.Range("B9").Value = 90000000
$ws.Range("G9").Value = 89900000

# --- Formatting updates ---------------------------------------------------
# Row 9 (summary totals row): bigger font + wrap text on
$row9 = $ws.Range("A9:J9")
$row9.Font.Size = 12
$row9.WrapText = $true

# Signature block (rows 12-19, columns C-H): bigger font only
$sigBlock = $ws.Range("C12:H19")
$sigBlock.Font.Size = 12
